# "add test case design"
# The F (实际值) and G (测试结果) columns for the data rows (rows 2-12) held
# placeholder/leftover test-result values that are removed here, leaving
# only the header labels in F1/G1. Column widths for E/F are re-asserted
# (23 and 22 characters) and the active selection is moved to F17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the previously-filled "actual value" / "test result" columns
# for every data row, leaving the header row (row 1) intact.
$ws.Range("F2:G12").ClearContents()

# Re-apply the column widths for columns E (5) and F (6) so they keep
# their custom width (23 and 22 characters respectively).
$ws.Columns.Item(5).ColumnWidth = 22.28
$ws.Columns.Item(6).ColumnWidth = 21.28

# Move the active selection like it was left after editing.
$ws.Range("F17").Select()
